$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (the "MuSCs" / "Resolving-Mac" target cluster row) entirely
$ws.Rows.Item(5).Delete()

# Update row 2 values
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.342322
$ws.Range("H2").Value = 1.026966
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.414593
$ws.Range("N2").Value = 4.243779
$ws.Range("O2").Value = 0.3478871232761722
$ws.Range("P2").Value = 0.3478871232761722
$ws.Range("Q2").Value = 0.484246304946
$ws.Range("R2").Value = 4.358216744513999
$ws.Range("S2").Value = 0.3478871232761722
$ws.Range("T2").Value = 0.3478871232761722

# Update row 3 values
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.342322
$ws.Range("H3").Value = 1.026966
$ws.Range("O3").Value = 0.5748520910875596
$ws.Range("P3").Value = 0.5748520910875596
$ws.Range("Q3").Value = 0.8001733389213332
$ws.Range("R3").Value = 7.201560050291998
$ws.Range("S3").Value = 0.5748520910875596
$ws.Range("T3").Value = 0.5748520910875596

# Update row 4 values
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.342322
$ws.Range("H4").Value = 1.026966
$ws.Range("M4").Value = 0.314161
$ws.Range("N4").Value = 0.942483
$ws.Range("O4").Value = 0.07726078563626818
$ws.Range("P4").Value = 0.07726078563626819
$ws.Range("Q4").Value = 0.107544221842
$ws.Range("R4").Value = 0.9678979965779998
$ws.Range("S4").Value = 0.07726078563626818
$ws.Range("T4").Value = 0.07726078563626819

$wb.Save()
